$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DuDad")
$ws.Range("D2").Value = "732-5005-ND"
